# Replace the division-problem answers throughout the worksheet table.
# Each old string is unique in the document, so a simple MatchWholeWord-off
# Find/Replace against the whole document content is safe and order-independent.
$d = $word.ActiveDocument

$replacements = @(
    @("909÷8=113, 5", "494÷2=247, 0"),
    @("389÷5=77, 4",  "405÷3=135, 0"),
    @("337÷3=112, 1", "185÷4=46, 1"),
    @("348÷8=43, 4",  "789÷7=112, 5"),
    @("825÷4=206, 1", "155÷6=25, 5"),
    @("476÷3=158, 2", "516÷8=64, 4"),
    @("108÷8=13, 4",  "111÷8=13, 7"),
    @("119÷9=13, 2",  "555÷3=185, 0"),
    @("487÷5=97, 2",  "357÷7=51, 0"),
    @("175÷8=21, 7",  "741÷4=185, 1"),
    @("483÷9=53, 6",  "521÷4=130, 1"),
    @("745÷4=186, 1", "769÷8=96, 1"),
    @("885÷4=221, 1", "761÷4=190, 1"),
    @("278÷3=92, 2",  "425÷3=141, 2"),
    @("134÷2=67, 0",  "370÷4=92, 2"),
    @("137÷5=27, 2",  "913÷2=456, 1"),
    @("316÷2=158, 0", "291÷2=145, 1"),
    @("192÷5=38, 2",  "658÷8=82, 2"),
    @("682÷3=227, 1", "131÷4=32, 3"),
    @("344÷3=114, 2", "592÷3=197, 1"),
    @("737÷4=184, 1", "811÷5=162, 1"),
    @("186÷3=62, 0",  "329÷5=65, 4"),
    @("215÷8=26, 7",  "518÷8=64, 6"),
    @("675÷7=96, 3",  "415÷9=46, 1"),
    @("820÷3=273, 1", "829÷3=276, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: could not find text to replace: $old"
    }
}
